$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsAssets = $wb.Worksheets.Item("Assets")

# --- Settings sheet -------------------------------------------------
# strSAPExecutable (row 7), strFolderToMove (row 10) and
# strNotepadExecutable (row 11) move out to Assets (below) and are
# replaced in-place: strComponentRoot/strTCode/strFindBankRoutingNumber
# shift up into rows 7-9.
$wsSettings.Range("A7").Value = "strComponentRoot"
$wsSettings.Range("B7").Value = "\\LRRBTUIPFSP100\Profiles\Uipath_26\Desktop\EssilorLuxottica Projects\Components\"
$wsSettings.Range("A8").Value = "strTCode"
$wsSettings.Range("B8").Value = "FLB2"
$wsSettings.Range("A9").Value = "strFindBankRoutingNumber"
$wsSettings.Range("B9").Value = 111000012

# Old row 10 (strFolderToMove) and its hyperlink go away; the cell
# keeps its leftover hyperlink style but no value.
$wsSettings.Range("A10:B10").ClearContents()
$wsSettings.Hyperlinks.Delete()

# Old rows 11 (strNotepadExecutable) and 12 (strFindBankRoutingNumber)
# are now fully empty.
$wsSettings.Range("A11:B12").ClearContents()

# --- Assets sheet -----------------------------------------------------
# The three settings removed above are now managed as Orchestrator
# Assets instead.
$wsAssets.Range("A4").Value = "strFolderToMove"
$wsAssets.Range("B4").Value = "Ess.LUX-7.0_BankFilesFolderToMove"
$wsAssets.Range("A5").Value = "strSAPExecutable"
$wsAssets.Range("B5").Value = "Ess.ALL_SAPExecutable"
$wsAssets.Range("A6").Value = "strNotepadExecutable"
$wsAssets.Range("B6").Value = "Ess.ALL_NotepadExecutable"

# --- Selection / active sheet bookkeeping -----------------------------
# Settings ends up as the active/selected tab, so select it last.
[void]$wsAssets.Range("A9").Select()
[void]$wsSettings.Range("A15").Select()
